$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "23.955.00"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.88%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.624.69"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("E4").Value = "  -0.51%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "308.34"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.26%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9987"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.49%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3941"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.58%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3853"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.34%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.9973"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.362"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "49.67"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.18%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.08493"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.89%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "23.87"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -4.62%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.083"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.80%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.641"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.27%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.00001288"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.37%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "1.621.66"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.05%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "93.99"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.86%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06930"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.44%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "20.11"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -4.78%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.853"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.17%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.9978"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.62%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "13.47"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.42%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "23.952.00"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.82%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.481"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +5.68%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.844"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.35%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "22.24"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.81%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "156.81"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.31%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "140.76"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -3.04%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "5.309"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -7.89%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.900"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.51%  "
$ws.Range("E32").Value = "  -1.00%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.799.39"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.12%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.08156"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.72%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.9939"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.38%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.02911"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -3.27%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "6.648"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.93%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.2679"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.95%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.09163"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -3.03%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "10.37"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.99%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "13.72"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.89%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.433"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -4.20%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.7547"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -3.05%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "15.99"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.45%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.6933"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.09%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.478"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.01%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "4.075"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.64%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.9979"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.52%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.08258"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -3.46%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "136.03"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.31%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.212"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -6.79%  "
